# Penalty/Reward system data cleanup (unfinished attempt).
# Sheet "Weekly Quantity": drop the oldest week's row and the two
# newest rows added at the tail.
# Sheet "Monthly Trend": drop the oldest month's row and the one
# newest row added at the tail.

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Rows.Item(2).Delete()
$wsWeekly.Rows.Item(8).Delete()
$wsWeekly.Rows.Item(8).Delete()

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Rows.Item(2).Delete()
$wsMonthly.Rows.Item(6).Delete()
